$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2302.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2302.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6908.400000000001
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -7244.400000000001
$ws.Range("H74").Value = 7111
$ws.Range("I74").Value = 7624.875
$ws.Range("K74").Value = 7624.875
$ws.Range("M74").Value = -6688.875
$ws.Range("H77").Value = 7111
$ws.Range("I77").Value = 7624.875
$ws.Range("K77").Value = 38124.375
$ws.Range("M77").Value = -33444.375
$ws.Range("H80").Value = 1762.1154
$ws.Range("J80").Value = 1895.4706
$ws.Range("L80").Value = 5686.4118
$ws.Range("N80").Value = -7682.4118
$ws.Range("H83").Value = 1762.1154
$ws.Range("J83").Value = 1895.4706
$ws.Range("L83").Value = 17059.2354
$ws.Range("N83").Value = -27043.2354
$ws.Range("H125").Value = 3821.3635
$ws.Range("I125").Value = 3781.6667
$ws.Range("K125").Value = 34035.0003
$ws.Range("M125").Value = -31575.0003
$ws.Range("H133").Value = 71999
$ws.Range("J133").Value = 71999
$ws.Range("L133").Value = 71999
$ws.Range("N133").Value = -82119
$ws.Range("H134").Value = 89000
$ws.Range("J134").Value = 89000
$ws.Range("L134").Value = 89000
$ws.Range("N134").Value = -99140
$ws.Range("H136").Value = 135388.75
$ws.Range("J136").Value = 135388.75
$ws.Range("L136").Value = 135388.75
$ws.Range("N136").Value = -145588.75
$ws.Range("H137").Value = 52635268
$ws.Range("I137").Value = 111113350
$ws.Range("J137").Value = 4989.9
$ws.Range("K137").Value = 333340050
$ws.Range("L137").Value = 14969.7
$ws.Range("M137").Value = -333337500
$ws.Range("N137").Value = -20069.7
$ws.Range("H138").Value = 3231.0952
$ws.Range("I138").Value = 3022
$ws.Range("J138").Value = 3259.3513
$ws.Range("K138").Value = 9066
$ws.Range("L138").Value = 9778.053899999999
$ws.Range("M138").Value = -3926
$ws.Range("N138").Value = -20058.0539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2021.6923
$ws.Range("I45").Value = 2028.2
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2028.2
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1651.2
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 5888814
$ws.Range("J61").Value = 100000000
$ws.Range("L61").Value = 100000000
$ws.Range("N61").Value = -100000424
$ws.Range("H74").Value = 2532556.8
$ws.Range("I74").Value = 3089276.2
$ws.Range("K74").Value = 3089276.2
$ws.Range("M74").Value = -3088402.2
$ws.Range("H77").Value = 2532556.8
$ws.Range("I77").Value = 3089276.2
$ws.Range("K77").Value = 15446381
$ws.Range("M77").Value = -15442013
$ws.Range("H136").Value = 5888814
$ws.Range("J136").Value = 100000000
$ws.Range("L136").Value = 300000000
$ws.Range("N136").Value = -300005100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 67516
$ws.Range("J20").Value = 1322.6
$ws.Range("K20").Value = 67516
$ws.Range("L20").Value = 1322.6
$ws.Range("M20").Value = -67269
$ws.Range("N20").Value = -1816.6
$ws.Range("H81").Value = 66500
$ws.Range("J81").Value = 66500
$ws.Range("L81").Value = 66500
$ws.Range("N81").Value = -68622
$ws.Range("H82").Value = 24866.785
$ws.Range("I82").Value = 10454.375
$ws.Range("J82").Value = 44083.332
$ws.Range("K82").Value = 10454.375
$ws.Range("L82").Value = 44083.332
$ws.Range("M82").Value = -10071.375
$ws.Range("N82").Value = -44849.332
$ws.Range("H84").Value = 66500
$ws.Range("J84").Value = 66500
$ws.Range("L84").Value = 199500
$ws.Range("N84").Value = -210108
$ws.Range("H85").Value = 24866.785
$ws.Range("I85").Value = 10454.375
$ws.Range("J85").Value = 44083.332
$ws.Range("K85").Value = 10454.375
$ws.Range("L85").Value = 44083.332
$ws.Range("M85").Value = -9128.375
$ws.Range("N85").Value = -46735.332
$ws.Range("H134").Value = 5380767
$ws.Range("I134").Value = 4570.2666
$ws.Range("K134").Value = 13710.7998
$ws.Range("M134").Value = -11175.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3627499.8
$ws.Range("I58").Value = 3132.75
$ws.Range("J58").Value = 5560495.5
$ws.Range("K58").Value = 3132.75
$ws.Range("L58").Value = 5560495.5
$ws.Range("M58").Value = -2929.75
$ws.Range("N58").Value = -5560901.5
$ws.Range("H118").Value = 73324
$ws.Range("H136").Value = 3627499.8
$ws.Range("I136").Value = 3132.75
$ws.Range("J136").Value = 5560495.5
$ws.Range("K136").Value = 9398.25
$ws.Range("L136").Value = 16681486.5
$ws.Range("M136").Value = -6848.25
$ws.Range("N136").Value = -16686586.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 13291.667
$ws.Range("I82").Value = 6500
$ws.Range("K82").Value = 19500
$ws.Range("M82").Value = -19094
$ws.Range("H85").Value = 13291.667
$ws.Range("I85").Value = 6500
$ws.Range("K85").Value = 19500
$ws.Range("M85").Value = -18096
$ws.Range("H86").Value = 148
$ws.Range("J86").Value = 198
$ws.Range("L86").Value = 594
$ws.Range("N86").Value = -2966
$ws.Range("H89").Value = 148
$ws.Range("J89").Value = 198
$ws.Range("L89").Value = 1782
$ws.Range("N89").Value = -13638
$ws.Range("H113").Value = 1072.2307
$ws.Range("J113").Value = 1186.4
$ws.Range("L113").Value = 3559.2
$ws.Range("N113").Value = -7899.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 18000
$ws.Range("J46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("N46").Value = -18312
$ws.Range("H70").Value = 29551.28
$ws.Range("J70").Value = 63559.293
$ws.Range("L70").Value = 63559.293
$ws.Range("N70").Value = -64099.293
$ws.Range("H73").Value = 29551.28
$ws.Range("J73").Value = 63559.293
$ws.Range("L73").Value = 63559.293
$ws.Range("N73").Value = -65431.293

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6719.4
$ws.Range("I61").Value = 5025.129
$ws.Range("K61").Value = 5025.129
$ws.Range("M61").Value = -4823.129
$ws.Range("H113").Value = 6719.4
$ws.Range("I113").Value = 5025.129
$ws.Range("K113").Value = 5025.129
$ws.Range("M113").Value = -2855.129
$ws.Range("H132").Value = 7793317
$ws.Range("I132").Value = 23372506
$ws.Range("J132").Value = 3722.5
$ws.Range("K132").Value = 70117518
$ws.Range("L132").Value = 11167.5
$ws.Range("M132").Value = -70114988
$ws.Range("N132").Value = -16227.5
$ws.Range("H133").Value = 67163
$ws.Range("J133").Value = 67163
$ws.Range("L133").Value = 67163
$ws.Range("N133").Value = -72223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 104246.63
$ws.Range("I122").Value = 2083
$ws.Range("J122").Value = 226843
$ws.Range("K122").Value = 6249
$ws.Range("L122").Value = 680529
$ws.Range("M122").Value = -3799
$ws.Range("N122").Value = -685429
$ws.Range("H136").Value = 15219326
$ws.Range("I136").Value = 3346557
$ws.Range("J136").Value = 66667996
$ws.Range("K136").Value = 10039671
$ws.Range("L136").Value = 200003988
$ws.Range("M136").Value = -10037121
$ws.Range("N136").Value = -200009088
